$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cosmetic: the built-in "Normal" cell style was renamed (locale change).
try {
    $wb.Styles("Normal").Name = "Normalny"
} catch {
}

# Rows 9-15 previously had no "B" cell at all. Give them the same style as
# column A in their row (style index 3) before (optionally) filling values.
for ($r = 9; $r -le 15; $r++) {
    $ws.Range("A$r").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
}

# Rows 2-11: fill the full set of config values.
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("A$r").Value = "ETH-USD"
    $ws.Range("B$r").Value = "Indicators"
    $ws.Range("C$r").Value = 60
    $ws.Range("D$r").Value = 6
    $ws.Range("E$r").Value = 100
    $ws.Range("F$r").Value = "Binance"
    $ws.Range("G$r").Value = "1h"
}

$ws.Range("E7").Select() | Out-Null
